# Edit: insert 4 new data rows (weekly update of Ajo prices) at the top of
# the block that starts at row 1132, pushing the existing rows down by 4.
# New rows correspond to the most recent reporting date (45132).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 1132; everything below (old 1132-1209) shifts
# down to become rows 1136-1213.
$ws.Rows("1132:1135").Insert()

# Row 1132: Ajo, Chino, Primera
$ws.Range("A1132").Value = 6
$ws.Range("B1132").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1132").Value = "Metropolitana"
$ws.Range("D1132").Value = 45132
$ws.Range("E1132").Value = 13
$ws.Range("F1132").Value = 100112003
$ws.Range("G1132").Value = "Ajo"
$ws.Range("H1132").Value = "Chino"
$ws.Range("I1132").Value = "Primera"
$ws.Range("J1132").Value = 2200
$ws.Range("K1132").Value = 17000
$ws.Range("L1132").Value = 17500
$ws.Range("M1132").Value = 17159
$ws.Range("N1132").Value = "$/caja 10 kilos"
$ws.Range("O1132").Value = "China"
$ws.Range("P1132").Value = 1716
$ws.Range("Q1132").Value = 10
$ws.Range("R1132").Value = "Hortaliza"

# Row 1133: Ajo, Rosado, 1a (guarda)
$ws.Range("A1133").Value = 6
$ws.Range("B1133").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1133").Value = "Metropolitana"
$ws.Range("D1133").Value = 45132
$ws.Range("E1133").Value = 13
$ws.Range("F1133").Value = 100112003
$ws.Range("G1133").Value = "Ajo"
$ws.Range("H1133").Value = "Rosado"
$ws.Range("I1133").Value = "1a (guarda)"
$ws.Range("J1133").Value = 500
$ws.Range("K1133").Value = 25000
$ws.Range("L1133").Value = 25000
$ws.Range("M1133").Value = 25000
$ws.Range("N1133").Value = "$/trenza 50 unidades"
$ws.Range("O1133").Value = "Región Metropolitana"
$ws.Range("P1133").Value = 2500
$ws.Range("Q1133").Value = 10
$ws.Range("R1133").Value = "Hortaliza"

# Row 1134: Ajo, Rosado, 2a (guarda)
$ws.Range("A1134").Value = 6
$ws.Range("B1134").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1134").Value = "Metropolitana"
$ws.Range("D1134").Value = 45132
$ws.Range("E1134").Value = 13
$ws.Range("F1134").Value = 100112003
$ws.Range("G1134").Value = "Ajo"
$ws.Range("H1134").Value = "Rosado"
$ws.Range("I1134").Value = "2a (guarda)"
$ws.Range("J1134").Value = 300
$ws.Range("K1134").Value = 20000
$ws.Range("L1134").Value = 20000
$ws.Range("M1134").Value = 20000
$ws.Range("N1134").Value = "$/trenza 50 unidades"
$ws.Range("O1134").Value = "Región Metropolitana"
$ws.Range("P1134").Value = 2000
$ws.Range("Q1134").Value = 10
$ws.Range("R1134").Value = "Hortaliza"

# Row 1135: Ajo, Rosado, 3a (guarda)
$ws.Range("A1135").Value = 6
$ws.Range("B1135").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1135").Value = "Metropolitana"
$ws.Range("D1135").Value = 45132
$ws.Range("E1135").Value = 13
$ws.Range("F1135").Value = 100112003
$ws.Range("G1135").Value = "Ajo"
$ws.Range("H1135").Value = "Rosado"
$ws.Range("I1135").Value = "3a (guarda)"
$ws.Range("J1135").Value = 200
$ws.Range("K1135").Value = 13000
$ws.Range("L1135").Value = 13000
$ws.Range("M1135").Value = 13000
$ws.Range("N1135").Value = "$/trenza 50 unidades"
$ws.Range("O1135").Value = "Región Metropolitana"
$ws.Range("P1135").Value = 1300
$ws.Range("Q1135").Value = 10
$ws.Range("R1135").Value = "Hortaliza"
